$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the contents of columns AB:AK and AM for data rows 2 through 32.
# Column AL ("PREVIOUS ACCOMPLISHMENT") is intentionally left untouched.
$lastRow = 32
$range = $ws.Range("AB2:AK$lastRow")
$range.ClearContents()

$range2 = $ws.Range("AM2:AM$lastRow")
$range2.ClearContents()
